# Insertion The_Saison et The_MoisSaison
# Marks the already-completed sub-steps of "Backoffice (page admin)" /
# "configSaison.html" as done ("ok") and records the extra (2h) estimate.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "configSaison.html" -> "configSaison.html (2h)"
$ws.Range("B153").Value = "configSaison.html (2h)"

# New "ok" markers for the freshly implemented sub-tasks.
$ws.Range("B150").Value = "ok"
$ws.Range("B151").Value = "ok"
$ws.Range("B155").Value = "ok"
$ws.Range("B156").Value = "ok"
$ws.Range("B157").Value = "ok"
$ws.Range("B158").Value = "ok"
$ws.Range("B160").Value = "ok"
$ws.Range("B163").Value = "ok"

# Re-affirm the bold styling of the "Backoffice (page admin)" block so the
# workbook's style table collapses the duplicate bold xf entry.
$ws.Range("A147").Font.Bold = $true
$ws.Range("A148").Font.Bold = $true
$ws.Range("A149").Font.Bold = $true
$ws.Range("A150").Font.Bold = $true

# Move the selection / scroll position to where the edits were made.
$null = $ws.Range("B154").Select()
